$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("model_summaries")

$ws.Range("N2").Value = 1.541175203748385
$ws.Range("P2").Value = 1.167852480108754
$ws.Range("Q2").Value = 0.9580622735172457
$ws.Range("R2").Value = 3.232744463347607
$ws.Range("T2").Value = 2.654117159539461
$ws.Range("U2").Value = 0.842716400745953
